# Leads template update:
#  - add "Submission Date" column (F)
#  - replace sample leads with Ahmad / Zein rows incl. mailto hyperlinks on the Email column
#  - drop the third sample row (Mike Davis)
#  - widen the Lead Description column and size the new Submission Date column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "Phone Number"
$ws.Range("D1").Value = "Subject"
$ws.Range("E1").Value = "Lead Description"

# --- Phone / subject / description stay as the existing sample copy ---
$ws.Range("C2").Value = "555-123-4567"
$ws.Range("D2").Value = "Interested in flooring services"
$ws.Range("E2").Value = "Looking for hardwood flooring installation for living room and kitchen. Budget around `$5000."

$ws.Range("C3").Value = "555-987-6543"
$ws.Range("D3").Value = "Commercial flooring inquiry"
$ws.Range("E3").Value = "Need commercial grade flooring for office space, approximately 2000 sq ft."

# --- New "Submission Date" column ---
$ws.Range("F1").Value = "Submission Date"

# --- New lead emails (with mailto hyperlinks) ---
$ws.Range("B2").Value = "abc@gmail.om"
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:abc@gmail.om", "", "", "abc@gmail.om") | Out-Null

$ws.Range("B3").Value = "xyz@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:xyz@gmail.com", "", "", "xyz@gmail.com") | Out-Null

# --- New lead names ---
$ws.Range("A2").Value = "Ahmad"
$ws.Range("A3").Value = "Zein"

# --- Submission dates ---
$d2 = Get-Date -Year 2025 -Month 7 -Day 11 -Hour 0 -Minute 0 -Second 0
$ws.Range("F2").Value = $d2
$ws.Range("F2").NumberFormat = "m/d/yyyy"

$d3 = Get-Date -Year 2025 -Month 7 -Day 2 -Hour 0 -Minute 0 -Second 0
$ws.Range("F3").Value = $d3
$ws.Range("F3").NumberFormat = "m/d/yyyy"

# --- Drop the old 4th sample row (Mike Davis) ---
$ws.Range("A4:F4").Delete() | Out-Null

# --- Column sizing ---
$ws.Columns.Item(5).ColumnWidth = 79.16666666666667   # stored width -> 80
$ws.Columns.Item(6).ColumnWidth = 21.5                # stored width -> ~22.3 (new Submission Date column)

# --- Selection parity with the authored workbook ---
$ws.Range("B12").Select() | Out-Null
